{"js": "// The edited report swaps the CNPJ value and zeroes out the\n// \"Revenda de mercadorias com dispensa de emiss\u00e3o de documento fiscal\"\n// revenue (line I), along with its two roll-up totals (lines III and X),\n// which previously all echoed the same 2000.00 figure.\nconst body = context.document.body;\n\n// 1) CNPJ: 1106462000199 -> 312\nconst cnpjResults = body.search(\"1106462000199\", { matchCase: true });\ncnpjResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < cnpjResults.items.length; i++) {\n  cnpjResults.items[i].insertText(\"312\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Every occurrence of the revenue figure 2000.00 -> 0.00\n// (line I, line III total, and line X grand total).\nconst amountResults = body.search(\"2000.00\", { matchCase: true });\namountResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < amountResults.items.length; i++) {\n  amountResults.items[i].insertText(\"0.00\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The edited report swaps the CNPJ value and zeroes out the\n# \"Revenda de mercadorias com dispensa de emiss\u00e3o de documento fiscal\"\n# revenue (line I), along with its two roll-up totals (lines III and X),\n# which previously all echoed the same 2000.00 figure.\n$d = $word.ActiveDocument\n\n# 1) CNPJ: 1106462000199 -> 312\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"1106462000199\"\n$find.Replacement.Text = \"312\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Every occurrence of the revenue figure 2000.00 -> 0.00\n# (line I, line III total, and line X grand total).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"2000.00\"\n$find2.Replacement.Text = \"0.00\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
